# Scen_B_SYS_MaxGrowthRates.xlsx -- add max growth constraints for
# ambient heat and district heat in the SRV sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRV")

# ---------------------------------------------------------------------
# 1) Make room in the lookup table (rows 14-18) for two new data rows
#    (AmbientHeat, DistrictHeat) and a blank spacer + a repeated
#    "Max growth rate / Starting value" header row, matching the new
#    layout:
#       17: header (Max growth rate | Starting value (GW))
#       18: SRV-PU / HPs
#       19: SRV-CS / HPs
#       20: blank spacer (styled)
#       21: header (Max growth rate | Starting value (PJ))
#       22: SRV / Biomass
#       23: SRV / Biogas
#       24: SRV / AmbientHeat   (new)
#       25: SRV / DistrictHeat  (new)
# ---------------------------------------------------------------------

# Insert 2 rows before old row 14 -> becomes the new UC formula rows 13/14
$ws.Rows.Item(13).Resize(2).Insert()
# Insert 1 more row before the (now shifted) header row so it lands on 17
$ws.Rows.Item(16).Insert()
# Insert 2 rows (spacer + header) before the (now shifted) Biomass row
$ws.Rows.Item(20).Resize(2).Insert()

# ---------------------------------------------------------------------
# 2) Fix up the two repeated-header rows and the re-numbered lookup rows
# ---------------------------------------------------------------------

# Row 17 header (was row 14) -- unchanged content, just confirming values
$ws.Range("C17").Value = "Max growth rate"
$ws.Range("D17").Value = "Starting value (GW)"

# Row 21 header for the PJ-valued rows (Biomass/Biogas/AmbientHeat/DistrictHeat)
$ws.Range("C21").Value = "Max growth rate"
$ws.Range("D21").Value = "Starting value (PJ)"

# Blank spacer row 20, C20 carries the numeric style used elsewhere in the table
$ws.Range("C20").Value = ""

# New lookup rows 24 (AmbientHeat) and 25 (DistrictHeat)
$ws.Range("A24").Value = "SRV"
$ws.Range("B24").Value = "AmbientHeat"
$ws.Range("C24").Value = 0.05
$ws.Range("D24").Value = 0.3

$ws.Range("A25").Value = "SRV"
$ws.Range("B25").Value = "DistrictHeat"
$ws.Range("C25").Value = 0.05
$ws.Range("D25").Value = 0.3

# ---------------------------------------------------------------------
# 3) Apply the same number format/style used by the sibling cells so the
#    new rows look like the rest of the lookup table (s="11" on column C)
# ---------------------------------------------------------------------
$ws.Range("C18").Copy() | Out-Null
$ws.Range("C24:C25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Fill in the new UC constraint rows 13 (AmbientHeat) and 14 (DistrictHeat)
# ---------------------------------------------------------------------

# Row 13 - AmbientHeat
$ws.Range("B13").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A24,"MaxGrowth",B24)'
$ws.Range("C13").Value = "ACT, GROWTH"
$ws.Range("F13").Value = "SRVAHT"
$ws.Range("G13").Value = "FT*"
$ws.Range("H13").Value = 2021
$ws.Range("I13").Value = "LO"
$ws.Range("J13").Formula = "=1+C24"
$ws.Range("K13").Value = 1
$ws.Range("L13").Formula = "=-D24"
$ws.Range("M13").Value = 5
$ws.Range("N13").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A24, "maximum growth rate of",B24)'

# Row 14 - DistrictHeat
$ws.Range("B14").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A25,"MaxGrowth",B25)'
$ws.Range("C14").Value = "ACT, GROWTH"
$ws.Range("F14").Value = "SRVHET"
$ws.Range("G14").Value = "FT*"
$ws.Range("H14").Value = 2021
$ws.Range("I14").Value = "LO"
$ws.Range("J14").Formula = "=1+C25"
$ws.Range("K14").Value = 1
$ws.Range("L14").Formula = "=-D25"
$ws.Range("M14").Value = 5
$ws.Range("N14").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A25, "maximum growth rate of",B25)'

# Match styling (J/L columns use custom numeric styles elsewhere in the table)
$ws.Range("J12").Copy() | Out-Null
$ws.Range("J13:J14").PasteSpecial(-4122) | Out-Null
$ws.Range("L12").Copy() | Out-Null
$ws.Range("L13:L14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5) The RHS code used in the header row (L5/L10) changed from the
#    year-specific "UC_RHSRTS~2018" to the generic "UC_RHSRTS"
# ---------------------------------------------------------------------
$ws.Range("L5").Value = "UC_RHSRTS"
$ws.Range("L10").Value = "UC_RHSRTS"

# ---------------------------------------------------------------------
# 6) Column width tweaks for the now-wider L/M/N area (L narrower, new
#    M column sized, N widened to fit the longer constraint names)
# ---------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 11.328125
$ws.Columns.Item(13).ColumnWidth = 13.671875
$ws.Columns.Item(14).ColumnWidth = 39.16015625

# ---------------------------------------------------------------------
# 7) Active-sheet/selection bookkeeping: the author ended their session
#    with SRV active (cell M29 selected) instead of SUP (cell J14).
# ---------------------------------------------------------------------
$sup = $wb.Worksheets.Item("SUP")
$sup.Activate()
$sup.Range("J14").Select()

$ws.Activate()
$ws.Range("M29").Select()
